$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: Sroufe & Gopalakrishna-Remani, SAGE / Organization and Environment
$ws.Range("A10").Value = 2019
$ws.Range("B10").Value = "SAGE"
$ws.Range("C10").Value = "Organization and Environment"
$ws.Range("D10").Value = "Management, Social Sustainability, Reputation, and Financial Performance Relationships: An Empirical Examination of U.S. Firms"
$ws.Range("E10").Value = "Robert Sroufe, Venugopal Gopalakrishna-Remani"
$ws.Range("F10").Value = "return on asset, return on investment, net profit margin"
$ws.Range("G10").Value = 3
$ws.Range("H10").Value = "Dependen"
$ws.Range("I10").Value = "Amerika"
$ws.Range("J10").Value = "PLS-SEM"
$ws.Range("K10").Value = "SmartPLS"
$ws.Range("L10").Value = "Fortune 500 firms simultaneously listed in the Newsweek Green rankings, The Corporate Knights Global 100, and the 100 Best Corporate Citizens lists"

# Row 11: Sujit & Rajesh, SAGE / SAGE Open
$ws.Range("A11").Value = 2016
$ws.Range("B11").Value = "SAGE"
$ws.Range("C11").Value = "SAGE Open"
$ws.Range("D11").Value = "Determinants of Discretionary Investments: Evidence From Indian Food Industry"
$ws.Range("E11").Value = "K. S. Sujit, B. K. Rajesh"
$ws.Range("F11").Value = "Profit before interest, tax, depreciation and amortization divided by total income; Net income/sales; Net income divided by total net worth; Net income divided by total assets; Profit after tax divided by total assets; Profit after tax divided by net worth; Profit after tax divided by total income"
$ws.Range("G11").Value = 7
$ws.Range("H11").Value = "Dependen"
$ws.Range("I11").Value = "India"
$ws.Range("J11").Value = "PLS-SEM"
$ws.Range("K11").Value = "SmartPLS"
$ws.Range("L11").Value = "Indian Food Industry"

# Two extra blank cells per new row (M, N) match the extended used range
$ws.Range("M10").Value = $null
$ws.Range("N10").Value = $null
$ws.Range("M11").Value = $null
$ws.Range("N11").Value = $null

# Apply the same center/center (no wrap) alignment style used elsewhere in the table
$range = $ws.Range("A10:N11")
$range.HorizontalAlignment = -4108
$range.VerticalAlignment = -4108
$range.WrapText = $false

$ws.Range("I15").Select()
